# Refatoramento teste de cadastro
$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Planilha1")
$ws2 = $wb.Worksheets.Item("Planilha2")

# Atualiza o usuario/email usados na massa de teste de cadastro (Planilha1)
$ws1.Range("A2").Value = "bvilhena"
$ws1.Range("B2").Value = "rodrigovil@rstinet.com"

# Mantem a selecao da Planilha2 e ajusta a selecao/aba ativa da Planilha1
$ws2.Range("E7").Select()
$ws1.Range("D7").Select()
$ws1.Activate()

$wb.Windows.Item(1).Width = 24240
